$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.174.03'
$ws.Range("E2").Value = '  -3.50%  '
$ws.Range("D3").Value = '3.505.20'
$ws.Range("E3").Value = '  -4.89%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.42%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.88%  '
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("D8").Value = '3.498.70'
$ws.Range("E8").Value = '  -4.90%  '
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("E10").Value = '  -6.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.71'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.83%  '
$ws.Range("E12").Value = '  -2.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '47.11'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.02%  '
$ws.Range("E14").Value = '  -3.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '675.18'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.27%  '
$ws.Range("D16").Value = '4.066.41'
$ws.Range("E16").Value = '  -5.01%  '
$ws.Range("E17").Value = '  -3.40%  '
$ws.Range("D18").Value = '69.133.66'
$ws.Range("E18").Value = '  -3.70%  '
$ws.Range("D19").Value = '3.506.70'
$ws.Range("E19").Value = '  -4.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.46'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.18'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.24%  '
$ws.Range("E23").Value = '  -4.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -9.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '97.94'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.80%  '
$ws.Range("E26").Value = '  -4.52%  '
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  -6.89%  '
$ws.Range("E30").Value = '  -7.30%  '
$ws.Range("E31").Value = '  -7.29%  '
$ws.Range("E32").Value = '  -5.95%  '
$ws.Range("E33").Value = '  -8.11%  '
$ws.Range("B34").Value = 'Mantle'
$ws.Range("C34").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.36'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.82%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.28'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '596.94'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.15%  '
$ws.Range("E37").Value = '  -15.22%  '
$ws.Range("E38").Value = '  -3.73%  '
$ws.Range("E39").Value = '  -4.08%  '
$ws.Range("E40").Value = '  -3.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("E42").Value = '  -5.88%  '
$ws.Range("E43").Value = '  -4.69%  '
$ws.Range("E44").Value = '  -7.55%  '
$ws.Range("D45").Value = '3.419.30'
$ws.Range("E45").Value = '  -10.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '33.41'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.26%  '
$ws.Range("E47").Value = '  -9.17%  '
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("E49").Value = '  -7.32%  '
$ws.Range("E50").Value = '  -0.78%  '
$ws.Range("E51").Value = '  +18.45%  '
